$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number and week-covering dates) ---
$ws.Range("A8").Value = "Volume 30   Number  1"
$ws.Range("C9").Value = "Report Covering the Week  1/2/2023  Through  1/8/2023"

# --- Crime Complaints data table updates (rows 14-30) ---
# Row 14
$ws.Range("D14").Value = "0"
$ws.Range("E14").Value = "***.*"
$ws.Range("M14").Value = -100
$ws.Range("N14").Value = -100

# Row 15
$ws.Range("D15").Value = "0"
$ws.Range("E15").Value = "***.*"
$ws.Range("F15").Value = 3
$ws.Range("H15").Value = 50
$ws.Range("I15").Value = 1
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = -66.666666666666

# Row 16
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 150
$ws.Range("F16").Value = 19
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = 111.111111111111
$ws.Range("I16").Value = 7
$ws.Range("J16").Value = 2
$ws.Range("K16").Value = 250
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 16.666666666666
$ws.Range("N16").Value = -85.416666666666

# Row 17
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 17
$ws.Range("E17").Value = -76.470588235294
$ws.Range("F17").Value = 30
$ws.Range("G17").Value = 41
$ws.Range("H17").Value = -26.829268292682
$ws.Range("I17").Value = 6
$ws.Range("J17").Value = 20
$ws.Range("K17").Value = -70
$ws.Range("L17").Value = -40
$ws.Range("M17").Value = 20
$ws.Range("N17").Value = -80.645161290322

# Row 18
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -42.857142857142
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = 14.285714285714
$ws.Range("I18").Value = 4
$ws.Range("J18").Value = 7
$ws.Range("K18").Value = -42.857142857142
$ws.Range("L18").Value = -42.857142857142
$ws.Range("M18").Value = -75
$ws.Range("N18").Value = -89.743589743589

# Row 19
$ws.Range("C19").Value = 22
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = 15.78947368421
$ws.Range("F19").Value = 68
$ws.Range("G19").Value = 88
$ws.Range("H19").Value = -22.727272727272
$ws.Range("I19").Value = 25
$ws.Range("J19").Value = 20
$ws.Range("K19").Value = 25
$ws.Range("L19").Value = 56.25
$ws.Range("M19").Value = 108.333333333333
$ws.Range("N19").Value = 25

# Row 20
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 30
$ws.Range("H20").Value = 233.333333333333
$ws.Range("I20").Value = 6
$ws.Range("J20").Value = 4
$ws.Range("K20").Value = 50
$ws.Range("L20").Value = 20
$ws.Range("M20").Value = 20
$ws.Range("N20").Value = -92

# Row 21
$ws.Range("C21").Value = 42
$ws.Range("D21").Value = 49
$ws.Range("E21").Value = -14.285714285714
$ws.Range("F21").Value = 166
$ws.Range("G21").Value = 165
$ws.Range("H21").Value = 0.60606060606
$ws.Range("I21").Value = 49
$ws.Range("J21").Value = 54
$ws.Range("K21").Value = -9.259259259259
$ws.Range("L21").Value = 6.521739130434
$ws.Range("M21").Value = 8.888888888888
$ws.Range("N21").Value = -77.419354838709

# Row 22
$ws.Range("C22").Value = "0"
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -100
$ws.Range("J22").Value = 2
$ws.Range("K22").Value = -100
$ws.Range("L22").Value = -100

# Row 23
$ws.Range("C23").Value = 4
$ws.Range("E23").Value = -20
$ws.Range("F23").Value = 17
$ws.Range("G23").Value = 18
$ws.Range("H23").Value = -5.555555555555
$ws.Range("I23").Value = 4
$ws.Range("J23").Value = 5
$ws.Range("K23").Value = -20
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 33.333333333333

# Row 24
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 76
$ws.Range("E24").Value = -60.526315789473
$ws.Range("F24").Value = 160
$ws.Range("G24").Value = 194
$ws.Range("H24").Value = -17.525773195876
$ws.Range("I24").Value = 30
$ws.Range("J24").Value = 76
$ws.Range("K24").Value = -60.526315789473
$ws.Range("L24").Value = -37.5
$ws.Range("M24").Value = 25

# Row 25
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 40
$ws.Range("F25").Value = 60
$ws.Range("G25").Value = 76
$ws.Range("H25").Value = -21.052631578947
$ws.Range("I25").Value = 16
$ws.Range("J25").Value = 13
$ws.Range("K25").Value = 23.076923076923
$ws.Range("L25").Value = 128.571428571429
$ws.Range("M25").Value = -20

# Row 26
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = "0"
$ws.Range("E26").Value = "***.*"
$ws.Range("F26").Value = 4
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 33.333333333333
$ws.Range("I26").Value = 1
$ws.Range("L26").Value = 0

# Row 27
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 100
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -20
$ws.Range("I27").Value = 2
$ws.Range("J27").Value = 1
$ws.Range("K27").Value = 100

# Row 28
$ws.Range("M28").Value = -100
$ws.Range("N28").Value = -100

# Row 29
$ws.Range("M29").Value = -100
$ws.Range("N29").Value = -100

# Row 30
$ws.Range("F30").Value = "0"

# --- Historical Perspective table updates (rows 37-43) ---
# Row 37
$ws.Range("J37").Value = 26
$ws.Range("K37").Value = -21.212121212121
$ws.Range("L37").Value = -21.212121212121
$ws.Range("M37").Value = -42.222222222222
$ws.Range("N37").Value = -38.095238095238

# Row 38
$ws.Range("J38").Value = 265
$ws.Range("K38").Value = -54.778156996587
$ws.Range("L38").Value = -66.916354556804
$ws.Range("M38").Value = -84.530064214827
$ws.Range("N38").Value = -85.367200441744

# Row 39
$ws.Range("J39").Value = 489
$ws.Range("K39").Value = 27.012987012987
$ws.Range("L39").Value = 9.887640449438
$ws.Range("M39").Value = -34.625668449197
$ws.Range("N39").Value = -41.716328963051

# Row 43
$ws.Range("J43").Value = 2295
$ws.Range("K43").Value = -34.070669347888
$ws.Range("L43").Value = -46.189917936694
$ws.Range("M43").Value = -73.590333716916
$ws.Range("N43").Value = -78.203058220153
